$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (Developpeur) and column I (status) values added for rows 2-53,
# mirroring the new "Developpeur" / status tracking columns added to the task list.
$ws.Range("H2").Value = "Chef de projet + développeur"
$ws.Range("I2").Value = "ok"
$ws.Range("H3").Value = "Visual studio + netbeans"
$ws.Range("I3").Value = "ok"
$ws.Range("H4").Value = "Justin"
$ws.Range("I4").Value = "en cours"
$ws.Range("H5").Value = "Justin"
$ws.Range("I5").Value = "attente"
$ws.Range("H6").Value = "Justin"
$ws.Range("I6").Value = "attente"
$ws.Range("H7").Value = "Justin"
$ws.Range("I7").Value = "attente"
$ws.Range("H8").Value = "Justin"
$ws.Range("I8").Value = "attente"
$ws.Range("H9").Value = "Luka"
$ws.Range("H10").Value = "Luka"
$ws.Range("H11").Value = "Julien"
$ws.Range("I11").Value = "ok"
$ws.Range("H15").Value = "Julien"
$ws.Range("H16").Value = "Luka"
$ws.Range("H17").Value = "Olivier"
$ws.Range("H18").Value = "Luka"
$ws.Range("H19").Value = "Luka"
$ws.Range("H20").Value = "Luka / Julien"
$ws.Range("H21").Value = "Justin"
$ws.Range("H22").Value = "Justin"
$ws.Range("H23").Value = "Olivier"
$ws.Range("H24").Value = "Julien"
$ws.Range("H25").Value = "Julien"
$ws.Range("H26").Value = "Justin / Julien"
$ws.Range("H27").Value = "Luka / Julien"
$ws.Range("H28").Value = "Luka / Julien"
$ws.Range("H29").Value = "Olivier"
$ws.Range("H30").Value = "Olivier"
$ws.Range("H31").Value = "Olivier"
$ws.Range("H32").Value = "Olivier"
$ws.Range("H33").Value = "Olivier"
$ws.Range("H34").Value = "Julien"
$ws.Range("H35").Value = "Julien"
$ws.Range("H36").Value = "Luka"
$ws.Range("H37").Value = "Justin"
$ws.Range("H38").Value = "Justin"
$ws.Range("H39").Value = "Justin"
$ws.Range("H40").Value = "Olivier"
$ws.Range("H41").Value = "Justin"
$ws.Range("H42").Value = "Justin "
$ws.Range("H43").Value = "Luka"
$ws.Range("H44").Value = "Luka"
$ws.Range("H45").Value = "Luka"
$ws.Range("H46").Value = "Olivier + Luka"
$ws.Range("H47").Value = "Olivier + Luka"
$ws.Range("H48").Value = "Olivier + Luka"
$ws.Range("H49").Value = "Justin + Julien"
$ws.Range("H50").Value = "Justin + Julien"
$ws.Range("H51").Value = "Julien"
$ws.Range("H52").Value = "Justin"
$ws.Range("H53").Value = "Justin"

# Restore the active selection that was in place when the workbook was last saved.
[void]$ws.Range("I5").Select()

